# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on every sheet that
#    shows it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" columns (Overview cols E & F, zh-cn/de-de col C)
#    from ~17.22 chars to ~13.41 chars, matching the regenerated report's
#    auto-fit widths.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

if ($wsOverview.Range("E2").Value2 -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
